# Auto-generated: apply scheduled-runner price/profit refresh to leve tracker sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2732.0278
$ws.Range("I15").Value = 2732.0278
$ws.Range("K15").Value = 8196.0834
$ws.Range("M15").Value = -8027.0834

$ws.Range("H19").Value = 606.17645
$ws.Range("I19").Value = 341.83334
$ws.Range("J19").Value = 750.36365
$ws.Range("K19").Value = 341.83334
$ws.Range("L19").Value = 750.36365
$ws.Range("M19").Value = -166.83334
$ws.Range("N19").Value = -1100.36365

$ws.Range("H101").Value = 635.05884
$ws.Range("I101").Value = 398.5
$ws.Range("J101").Value = 1202.8
$ws.Range("K101").Value = 1195.5
$ws.Range("L101").Value = 3608.4
$ws.Range("M101").Value = 426.5
$ws.Range("N101").Value = -6852.4

$ws.Range("H104").Value = 741.5
$ws.Range("I104").Value = 888.6667
$ws.Range("K104").Value = 2666.0001
$ws.Range("M104").Value = -919.0001000000002

$ws.Range("H125").Value = 1461.4286
$ws.Range("I125").Value = 921.3333
$ws.Range("J125").Value = 1866.5
$ws.Range("K125").Value = 8291.9997
$ws.Range("L125").Value = 16798.5
$ws.Range("M125").Value = -5831.9997
$ws.Range("N125").Value = -21718.5

$ws.Range("H132").Value = 893191.4
$ws.Range("I132").Value = 1602.1428
$ws.Range("J132").Value = 2453472.5
$ws.Range("K132").Value = 4806.428400000001
$ws.Range("L132").Value = 7360417.5
$ws.Range("M132").Value = -2276.428400000001
$ws.Range("N132").Value = -7365477.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19920.459
$ws.Range("I32").Value = 18734.672
$ws.Range("K32").Value = 18734.672
$ws.Range("M32").Value = -18447.672

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3189.6553
$ws.Range("I31").Value = 3150.9092
$ws.Range("J31").Value = 3213.3333
$ws.Range("K31").Value = 3150.9092
$ws.Range("L31").Value = 3213.3333
$ws.Range("M31").Value = -2855.9092
$ws.Range("N31").Value = -3803.3333

$ws.Range("H34").Value = 3189.6553
$ws.Range("I34").Value = 3150.9092
$ws.Range("J34").Value = 3213.3333
$ws.Range("K34").Value = 3150.9092
$ws.Range("L34").Value = 3213.3333
$ws.Range("M34").Value = -2948.9092
$ws.Range("N34").Value = -3617.3333

$ws.Range("H107").Value = 668.2857
$ws.Range("I107").Value = 575.6
$ws.Range("K107").Value = 575.6
$ws.Range("M107").Value = 1344.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 2860
$ws.Range("I94").Value = 1200
$ws.Range("J94").Value = 3275
$ws.Range("K94").Value = 3600
$ws.Range("L94").Value = 9825
$ws.Range("M94").Value = -2924
$ws.Range("N94").Value = -11177

$ws.Range("H97").Value = 2701.4443
$ws.Range("I97").Value = 3243.5715
$ws.Range("J97").Value = 804
$ws.Range("K97").Value = 9730.7145
$ws.Range("L97").Value = 2412
$ws.Range("M97").Value = -9234.7145
$ws.Range("N97").Value = -3404

$ws.Range("H100").Value = 2510.5264
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 2538.889
$ws.Range("K100").Value = 6000
$ws.Range("L100").Value = 7616.667
$ws.Range("M100").Value = -5189
$ws.Range("N100").Value = -9238.667000000001

$ws.Range("H103").Value = 2776.9565
$ws.Range("J103").Value = 4023.7334
$ws.Range("L103").Value = 12071.2002
$ws.Range("N103").Value = -13829.2002

$ws.Range("H106").Value = 2940.1
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 2940.1
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 8820.299999999999
$ws.Range("N106").Value = -10712.3

$ws.Range("H109").Value = 2595.1538
$ws.Range("I109").Value = 789.5
$ws.Range("K109").Value = 2368.5
$ws.Range("M109").Value = -1328.5

$ws.Range("H112").Value = 14496192
$ws.Range("I112").Value = 2347.5
$ws.Range("J112").Value = 22226242
$ws.Range("K112").Value = 7042.5
$ws.Range("L112").Value = 66678726
$ws.Range("M112").Value = -5934.5
$ws.Range("N112").Value = -66680942

$ws.Range("H115").Value = 2929.2307
$ws.Range("J115").Value = 2923.3333
$ws.Range("L115").Value = 8769.999899999999
$ws.Range("N115").Value = -11119.9999

$ws.Range("H118").Value = 3959
$ws.Range("I118").Value = 449.66666
$ws.Range("J118").Value = 5275
$ws.Range("K118").Value = 1348.99998
$ws.Range("L118").Value = 15825
$ws.Range("M118").Value = -105.9999800000001
$ws.Range("N118").Value = -18311

$ws.Range("H121").Value = 95673864
$ws.Range("I121").Value = 1364.8334
$ws.Range("J121").Value = 139830400
$ws.Range("K121").Value = 4094.5002
$ws.Range("L121").Value = 419491200
$ws.Range("M121").Value = -2784.5002
$ws.Range("N121").Value = -419493820

$ws.Range("H122").Value = 851.1875
$ws.Range("J122").Value = 1083.909
$ws.Range("L122").Value = 9755.181
$ws.Range("N122").Value = -14655.181

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0

$ws.Range("H124").Value = 898.2

$ws.Range("H125").Value = 2722.7273

$ws.Range("H131").Value = 1062.6818
$ws.Range("J131").Value = 1428.5714
$ws.Range("L131").Value = 4285.7142
$ws.Range("N131").Value = -14365.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 70685.34
$ws.Range("I132").Value = 51492.5
$ws.Range("J132").Value = 113336.11
$ws.Range("K132").Value = 154477.5
$ws.Range("L132").Value = 340008.33
$ws.Range("M132").Value = -151947.5
$ws.Range("N132").Value = -345068.33

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1126.1428
$ws.Range("I93").Value = 1006.5455
$ws.Range("J93").Value = 1564.6666
$ws.Range("K93").Value = 1006.5455
$ws.Range("L93").Value = 1564.6666
$ws.Range("M93").Value = 241.4545000000001
$ws.Range("N93").Value = -4060.6666

$ws.Range("H122").Value = 3326.879
$ws.Range("I122").Value = 2920.7273
$ws.Range("J122").Value = 3529.9546
$ws.Range("K122").Value = 8762.1819
$ws.Range("L122").Value = 10589.8638
$ws.Range("M122").Value = -6312.1819
$ws.Range("N122").Value = -15489.8638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 42813.75
$ws.Range("I100").Value = 46192.91
$ws.Range("J100").Value = 39954.46
$ws.Range("K100").Value = 92385.82000000001
$ws.Range("L100").Value = 79908.92
$ws.Range("M100").Value = -91844.82000000001
$ws.Range("N100").Value = -80990.92

$ws.Range("H122").Value = 2454.9375
$ws.Range("J122").Value = 3897.5
$ws.Range("L122").Value = 11692.5
$ws.Range("N122").Value = -16592.5

# Cells removed by the refresh (values now blank, not zero)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M106").ClearContents()
$ws.Range("N123").ClearContents()
